# Auto-generated edit script: updates crypto price/volume table per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.171.48"
$ws.Range("E2").Value = "  +2.96%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.060.44"
$ws.Range("E3").Value = "  +2.63%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.37"
$ws.Range("E5").Value = "  +2.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  +2.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.47"
$ws.Range("E7").Value = "  +8.41%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  +3.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0813"
$ws.Range("E10").Value = "  +4.82%  "

$ws.Range("E11").Value = "  +2.57%  "

$ws.Range("E12").Value = "  +6.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.365.78"
$ws.Range("E13").Value = "  +2.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.29"
$ws.Range("E14").Value = "  +8.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.757"
$ws.Range("E15").Value = "  +3.52%  "

$ws.Range("E16").Value = "  +2.95%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.075.32"
$ws.Range("E17").Value = "  +4.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.007.86"
$ws.Range("E18").Value = "  +2.83%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.28"
$ws.Range("E19").Value = "  +1.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.97"
$ws.Range("E20").Value = "  +2.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0837"
$ws.Range("E21").Value = "  +3.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.29"
$ws.Range("E22").Value = "  +0.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("E25").Value = "  +5.03%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.46"
$ws.Range("E26").Value = "  +1.50%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.29"
$ws.Range("E27").Value = "  +4.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.133"
$ws.Range("E28").Value = "  +8.40%  "

$ws.Range("E29").Value = "  +2.70%  "

$ws.Range("E30").Value = "  +3.44%  "

$ws.Range("E31").Value = "  +3.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.57"
$ws.Range("E32").Value = "  +4.00%  "

$ws.Range("E33").Value = "  +3.29%  "

$ws.Range("E35").Value = "  +1.86%  "

$ws.Range("E36").Value = "  +1.85%  "

$ws.Range("E37").Value = "  +16.12%  "

$ws.Range("E38").Value = "  +6.12%  "

$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.534.73"
$ws.Range("E40").Value = "  +5.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.39"

$ws.Range("E42").Value = "  +3.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.99"
$ws.Range("E43").Value = "  +7.45%  "

$ws.Range("E44").Value = "  +4.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0925"
$ws.Range("E45").Value = "  +2.13%  "

$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.19"
$ws.Range("E46").Value = "  +1.81%  "

$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.13"
$ws.Range("E47").Value = "  +2.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.03"
$ws.Range("E48").Value = "  +3.44%  "

$ws.Range("E49").Value = "  +3.48%  "

$ws.Range("E50").Value = "  +0.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.252.21"
$ws.Range("E51").Value = "  +2.68%  "
